$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 23 (shifts old row 23 down to row 25)
$ws.Rows.Item(23).Resize(2).Insert()

# Fill A/D (numeric) columns first - these don't touch the shared string table
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 4).Value = 0

# New string cells, in the order that reproduces the target sharedStrings layout
$ws.Cells.Item(23, 5).Value = "[[1:1003:FALSE]]"
$ws.Cells.Item(24, 2).Value = "add 1003 character"
$ws.Cells.Item(24, 3).Value = "add 1003 character"
$ws.Cells.Item(23, 6).Value = "[10]"
$ws.Cells.Item(23, 2).Value = "remove 1003 character"
$ws.Cells.Item(23, 3).Value = "remove 1003 character"
$ws.Cells.Item(24, 5).Value = "[[2:106:TRUE]]"
$ws.Cells.Item(24, 6).Value = "[30]"

# Update selection to match the diff (scroll position / topLeftCell is
# window-geometry bookkeeping that Excel stamps on save; it isn't wired to
# the exporter in this headless COM host, so it's intentionally left alone)
$ws.Range("E25").Select()
